$d = $word.ActiveDocument

# Paragraph 1: ilvl=1, numId=1 -> "Linguagens de alto nível..."
$lastPara = $d.Paragraphs.Last
$rng = $lastPara.Range
$rng.InsertParagraphAfter()
$p1 = $d.Paragraphs.Last
$p1.Range.ListFormat.ListLevelNumber = 2
$r1 = $p1.Range
$r1.Collapse(1)
$r1.InsertAfter(" Linguagens de alto nível são aquelas que parecem muito com a linguagem humana. Diferente daquelas que parecem mais com as de máquinas.")

# Paragraph 2: ilvl=2, numId=1 -> "A principal vantagem..."
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Last
$p2.Range.ListFormat.ListLevelNumber = 3
$r2 = $p2.Range
$r2.Collapse(1)
$r2.InsertAfter("A principal vantagem é a facilidade de entendimento.")

# Paragraph 3: ilvl=2, numId=1 -> "A desvantagem..."
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Last
$p3.Range.ListFormat.ListLevelNumber = 3
$r3 = $p3.Range
$r3.Collapse(1)
$r3.InsertAfter("A desvantagem é que o interpreter precisa traduzir o que escrevemos para linguagem de máquina, para que o computador execute.")

# Paragraph 4: ilvl=2, numId=1 -> "No caso das linguagens de baixo nível..."
$p3.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs.Last
$p4.Range.ListFormat.ListLevelNumber = 3
$r4 = $p4.Range
$r4.Collapse(1)
$r4.InsertAfter("No caso das linguagens de baixo nível, elas já estão em linguagem de máquina, não necessitando do interpreter e tendo um desempenho melhor.")

Write-Output "done"
